# Update "想去人数" (interested-count) figures for several events in the
# "展览" and "全部类型" sheets, reflecting the newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 108
$ws1.Range("F4").Value = 1523
$ws1.Range("F5").Value = 219
$ws1.Range("F7").Value = 510
$ws1.Range("F8").Value = 9967
$ws1.Range("F14").Value = 6896
$ws1.Range("F15").Value = 1088
$ws1.Range("F16").Value = 640
$ws1.Range("F17").Value = 53
$ws1.Range("F18").Value = 203

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 108
$ws4.Range("F4").Value = 1523
$ws4.Range("F5").Value = 219
$ws4.Range("F8").Value = 511
$ws4.Range("F11").Value = 9967
$ws4.Range("F17").Value = 6896
$ws4.Range("F18").Value = 1088
$ws4.Range("F19").Value = 640
$ws4.Range("F20").Value = 53
$ws4.Range("F21").Value = 203
